# Fill in the previously-blank "B" column (most recent period) values on the
# VLO cash-flow sheet, and correct the existing Dividends Paid (Common) value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VLO")

$values = @{
    "B3"  = 2347000000.0
    "B4"  = -2480000000.0
    "B5"  = -3570000000.0
    "B7"  = 3146000000.0
    "B8"  = 16230000000.0
    "B9"  = 69000000.0
    "B10" = 899000000.0
    "B11" = 945000000.0
    "B12" = -1689000000.0
    "B14" = -47000000.0
    "B15" = -414000000.0
    "B16" = -2248000000.0
    "B17" = 3697000000.0
    "B18" = -23000000.0
    "B19" = -1599000000.0
    "B20" = -2803000000.0
    "B21" = 1833000000.0
    "B22" = 209000000.0
    "B23" = 739000000.0
    "B24" = 3313000000.0
    "B25" = 4052000000.0
    "B27" = -1807000000.0
    "B28" = 1053000000.0
    "B29" = -47000000.0
    "B30" = -23000000.0
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
